$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - preserve text type for numeric-looking strings
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "241.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.450"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05570"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.471"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.358"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8038"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.065"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1422"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07315"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03286"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02982"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09249"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001656"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.249"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04766"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005700"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006252"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001048"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003798"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001497"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0004174"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.971"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1303"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04221"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006953"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1043"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002904"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6788"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03076"
$ws.Range("D48").Style = "Normal"

# Coin / Link / Volume column updates (plain text, no conversion risk)
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "UpBots"
$ws.Range("C23").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("E23").Value = "22UpBotsUBXT"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
